# Revert the later "dump individual files" edit: point D/E (Abstract/Authors)
# for rows 2, 3, 4, 8 and 10 back at the earlier shared-string text (the one
# without the extra padding/prefix introduced by the reverted commit).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = 'id="Par1">As COVID-19 delivered a devastating blow around Europe, unprecedented mental health issues have arisen.

 This study outlines the Greek translation and validation of the Fear of COVID-19 Scale (FCV-19S) in the general population.

 A sample of 3029 Greek-speaking individuals completed the FCV-19S, as part of a large online survey reporting on the aftermath of the coronavirus on the psychological health of Greeks.

 A total of 2970 participants completed all items without missing values, and their responses were included in the statistical analysis.

 The confirmatory factor analysis suggested an acceptable model fit.

 The internal consistency measured by Cronbach’s alpha = .

87 for the whole scale.

 All hypothesized correlations were as expected, confirming the construct validity of the scale.

 The results of the analyses suggest that the psychometric properties of the Greek FCV-19S are sufficient.

'
$ws.Range("E2").Value = '[V.%Tsipropoulou%NULL%1,   V. A.%Nikopoulou%NULL%1,   Vasiliki%Holeva%vholeva@yahoo.gr%1,   Z.%Nasika%NULL%2,   Z.%Nasika%NULL%0,   I.%Diakogiannis%NULL%1,   S.%Sakka%NULL%1,   S.%Kostikidou%NULL%1,   C.%Varvara%NULL%1,   E.%Spyridopoulou%NULL%1,   E.%Parlapani%NULL%1]'
$ws.Range("D3").Value = 'id="Par1">The COVID-19 pandemic has impacted college students’ coursework, stress levels, and perceived health.

 Various estimates indicate that high proportions of college students have experienced increased amounts of stress (Dziech, Inside Higher Education, 2020; Hartocollis, New York Times, 2020), yet other sources indicate that many college students respond to these changes with resilience (Kelley, Cornell Chronicle, 2020).

 A method for assessing student anxiety regarding the pandemic is thus needed.

 The Fear of COVID-19 Scale (FCV-19S) (Ahorsu et al.

, International Journal of Mental Health and Addiction, 2020) is a seven-item scale which has been validated and shown to possess good psychometric qualities in studies of participants from multiple countries.

 The current study used a cross-sectional convenience sample of US college student participants (n = 237) and found that the FCV-19S has high reliability and validity as demonstrated by its internal consistency and strong one-factor solution.

 Scores on the FCV-19S were positively correlated with anxiety for students who were married or of Asian descent.

 Additionally, the FCV-19S was only moderately correlated with the Generalized Anxiety Disorder-7 Scale (GAD-7), suggesting that the FCV-19S may bring added utility to research and clinical practice with populations impacted by the COVID-19 pandemic.

'
$ws.Range("E3").Value = '[Catherine A.%Perz%perzc@uhv.edu%1,   Brent A.%Lang%NULL%2,   Brent A.%Lang%NULL%0,   Rick%Harrington%NULL%1]'
$ws.Range("D4").Value = 'id="Par1">Fear of novel coronavirus 2019 (COVID-19) may result in psychological health problems among different populations.

 Moreover, believing COVID-19 information and preventive COVID-19 infection behaviors are relevant constructs associated with fear of COVID-19. Therefore, the present study validated three instruments assessing fear, beliefs, and preventive behaviors related to COVID-19 among individuals with mental illness.

 Moreover, relationships between the three constructs were examined.

 Individuals with mental illness (N = 400; 178 females; mean age = 46.91 years) completed the Fear of COVID-19 Scale (FCV-19S), Believing COVID-19 Information Scale (BCIS), Preventive COVID-19 Infection Behaviors Scale (PCIBS), and Depression Anxiety Stress Scale-21 (DASS-21).

 The FCV-19S, BCIS, and PCIBS demonstrated a single-factor structure with satisfactory fit indices.

 Moreover, believing COVID-19 information positively and significantly associated with fear of COVID-19, and fear of COVID-19 negatively and significantly associated with preventive behaviors and positively and significantly associated with psychological distress.

 The FCV-19S, BCIS, and PCIBS may assist healthcare providers in assessing COVID-19-related information among individuals with mental illness.

 Consequently, relevant programs may be designed to help individuals with mental illness going through the period of crisis.

'
$ws.Range("E4").Value = '[Kun-Chia%Chang%NULL%1,   Wen-Li%Hou%NULL%1,   Amir H.%Pakpour%NULL%4,   Chung-Ying%Lin%cylin36933@gmail.com%3,   Chung-Ying%Lin%cylin36933@gmail.com%0,   Mark D.%Griffiths%NULL%0,   Mark D.%Griffiths%NULL%0]'
$ws.Range("D8").Value = 'id="Par1">Fear is a central emotional response to imminent threats such as the coronavirus-19 disease (COVID-19).

 The Fear of COVID-19 Scale (FCV-19S) assesses the severity of fear towards COVID-19. The present study examined the psychometric properties of the Arabic version of the FCV-19S.

 Using a forward-backward translation, the FCV-19S was translated into Arabic.

 An online survey using the Arabic versions of FCV-19S and the Hospital Anxiety and Depression Scale (HADS) was administered.

 Reliability and concurrent and confirmatory validity were examined.

 The dataset consisted of 693 Saudi participants.

 The internal consistency of the Arabic FCV-19S was satisfactory (α = .

88), with sound concurrent validity indicated by significant and positive correlations with HADS (r = .

66).

 The unidimensional structure of the FCV-19S was confirmed.

 The Arabic version of the FCV-19S is psychometrically robust and can be used in research assessing the psychological impact of COVID-19 among a Saudi adult population.

'
$ws.Range("E8").Value = '[Mohsen%Alyami%maly002@aucklanduni.ac.nz%1,   Marcus%Henning%NULL%2,   Marcus%Henning%NULL%0,   Christian U.%Krägeloh%NULL%2,   Christian U.%Krägeloh%NULL%0,   Hussain%Alyami%NULL%2,   Hussain%Alyami%NULL%0]'
$ws.Range("D10").Value = 'id="Par1">The recently developed Fear of COVID-19 Scale (FCV-19S) is a seven-item uni-dimensional scale that assesses the severity of fears of COVID-19. Given the rapid increase of COVID-19 cases in Bangladesh, we aimed to translate and validate the FCV-19S in Bangla.

 The forward-backward translation method was used to translate the English version of the questionnaire into Bangla.

 The reliability and validity properties of the Bangla FCV-19S were rigorously psychometrically evaluated (utilizing both confirmatory factor analysis and Rasch analysis) in relation to socio-demographic variables, national lockdown variables, and response to the Bangla Health Patient Questionnaire.

 The sample comprised 8550 Bangladeshi participants.

 The Cronbach α value for the Bangla FCV-19S was 0.871 indicating very good internal reliability.

 The results of the confirmatory factor analysis showed that the uni-dimensional factor structure of the FCV-19S fitted well with the data.

 The FCV-19S was significantly correlated with the nine-item Bangla Patient Health Questionnaire (PHQ-90) (r = 0.406, p &lt; 0.001).

 FCV-19S scores were significantly associated with higher worries concerning lockdown.

 Measurement invariance of the FCV-19S showed no differences with respect to age or gender.

 The Bangla version of FCV-19S is a valid and reliable tool with robust psychometric properties which will be useful for researchers carrying out studies among the Bangla speaking population in assessing the psychological impact of fear from COVID-19 infection during this pandemic.

'
$ws.Range("E10").Value = '[Najmuj%Sakib%NULL%1,   A. K. M. Israfil%Bhuiyan%NULL%1,   Sahadat%Hossain%NULL%1,   Firoj%Al Mamun%NULL%1,   Ismail%Hosen%NULL%1,   Abu Hasnat%Abdullah%NULL%1,   Md. Abedin%Sarker%NULL%1,   Mohammad Sarif%Mohiuddin%NULL%1,   Istihak%Rayhan%NULL%1,   Moazzem%Hossain%NULL%1,   Md. Tajuddin%Sikder%NULL%1,   David%Gozal%NULL%1,   Mohammad%Muhit%NULL%1,   S. M. Shariful%Islam%NULL%1,   Mark D.%Griffiths%mark.griffiths@ntu.ac.uk%0,   Amir H.%Pakpour%apakpour@qums.ac.ir%0,   Amir H.%Pakpour%apakpour@qums.ac.ir%0,   Mohammed A.%Mamun%mamunphi46@gmail.com%2,   Mohammed A.%Mamun%mamunphi46@gmail.com%0]'

# Writing multi-line text via COM auto-expands row height; put it back to the
# sheet default so only the cell values themselves show up as changed.
$ws.Rows.Item(2).AutoFit()
$ws.Rows.Item(3).AutoFit()
$ws.Rows.Item(4).AutoFit()
$ws.Rows.Item(8).AutoFit()
$ws.Rows.Item(10).AutoFit()
